$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.396.15'
$ws.Range("E2").Value = '  +0.48%  '

$ws.Range("D3").Value = '1.848.60'
$ws.Range("E3").Value = '  -0.42%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.37'
$ws.Range("E5").Value = '  +0.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4672'
$ws.Range("E7").Value = '  -1.45%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2730'
$ws.Range("E8").Value = '  -0.57%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06293'
$ws.Range("E9").Value = '  -2.04%  '

$ws.Range("D10").Value = '1.850.85'
$ws.Range("E10").Value = '  -0.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07467'
$ws.Range("E11").Value = '  +0.49%  '

$ws.Range("E12").Value = '  +1.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.931'
$ws.Range("E13").Value = '  -0.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '83.86'
$ws.Range("E14").Value = '  -1.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6204'
$ws.Range("E15").Value = '  -1.85%  '

$ws.Range("D16").Value = '30.332.05'
$ws.Range("E16").Value = '  +0.37%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.71'
$ws.Range("E18").Value = '  +1.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007309'
$ws.Range("E19").Value = '  -0.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.37'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.912'
$ws.Range("E22").Value = '  -3.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.872'
$ws.Range("E23").Value = '  -2.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '166.40'
$ws.Range("E24").Value = '  -0.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.146'
$ws.Range("E25").Value = '  -0.96%  '

$ws.Range("E26").Value = '  +0.25%  '

$ws.Range("E27").Value = '  +0.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1022'
$ws.Range("E28").Value = '  -0.16%  '

$ws.Range("E29").Value = '  -0.38%  '

$ws.Range("E30").Value = '  -3.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.816'
$ws.Range("E31").Value = '  -2.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04874'
$ws.Range("E32").Value = '  -0.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.139'
$ws.Range("E33").Value = '  -0.68%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7000'
$ws.Range("E34").Value = '  -3.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.690'
$ws.Range("E35").Value = '  +0.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01921'
$ws.Range("E36").Value = '  +0.28%  '

$ws.Range("E37").Value = '  +1.33%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.8645'
$ws.Range("E38").Value = '  -3.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '105.83'
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.937'
$ws.Range("E40").Value = '  -1.88%  '

$ws.Range("E41").Value = '  +0.64%  '

$ws.Range("E42").Value = '  +0.00%  '

$ws.Range("E43").Value = '  -1.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.077'
$ws.Range("E44").Value = '  +0.67%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.26'
$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1205'
$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.558'
$ws.Range("E47").Value = '  -2.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '33.35'
$ws.Range("E48").Value = '  +1.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05541'
$ws.Range("E49").Value = '  -0.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.345'
$ws.Range("E50").Value = '  -3.98%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3647'
$ws.Range("E51").Value = '  -1.53%  '
